# HeadersToClaimsMiddleware test data update
# The test data moved from a "ScopeProperties"/claims-based test sheet to a
# "HeadersToClaims" test sheet: the separate "Claims" data row for each test
# case (A, B, C) is removed, and the "Headers" / "Expected" rows get new
# values reflecting the HeadersToClaims middleware behavior.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerValue = 'header*X-UserScope=ABC&header*X-Role=admin&header*X-User=moe@stooges.org'

# Test case A
$ws.Cells.Item(2, 5).Value = 'A'
$ws.Cells.Item(2, 6).Value = 'Headers'
$ws.Cells.Item(2, 7).Value = $headerValue

$ws.Cells.Item(3, 5).Value = 'A'
$ws.Cells.Item(3, 6).Value = 'Expected'
$ws.Cells.Item(3, 7).Value = '[{"Type":"X-UserScope","Value":"ABC"}]'

# Test case B
$ws.Cells.Item(4, 5).Value = 'B'
$ws.Cells.Item(4, 6).Value = 'Headers'
$ws.Cells.Item(4, 7).Value = $headerValue

$ws.Cells.Item(5, 5).Value = 'B'
$ws.Cells.Item(5, 6).Value = 'Expected'
$ws.Cells.Item(5, 7).Value = '[{"Type":"X-Role","Value":"admin"},{"Type":"X-User","Value":"moe@stooges.org"}]'

# Test case C
$ws.Cells.Item(6, 5).Value = 'C'
$ws.Cells.Item(6, 6).Value = 'Headers'
$ws.Cells.Item(6, 7).Value = $headerValue

$ws.Cells.Item(7, 5).Value = 'C'
$ws.Cells.Item(7, 6).Value = 'Expected'
$ws.Cells.Item(7, 7).Value = '[{"Type":"X-UserScope","Value":"ABC"},{"Type":"X-Role","Value":"admin"},{"Type":"X-User","Value":"moe@stooges.org"}]'

# The old sheet had 3 rows per test case (Claims/Headers/Expected) spanning
# rows 2-10; the new layout only has 2 rows per test case (Headers/Expected)
# spanning rows 2-7, so remove the now-unused trailing rows.
$ws.Rows("8:10").Delete()

# Match the saved selection/active cell from the edited workbook.
$ws.Range("G8").Select()
